# Generate Report for Handoff
# - Overview sheet: refresh the "Latest HO Xliff Generate Date" timestamp
#   for the rows that were just (re)generated (rows 4-7).
# - zh-cn / de-de sheets: those same rows had their handoff re-run, so:
#     * Priority moves from "low" to "ht"
#     * Latest Handoff Datetime is refreshed to the new handoff time

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G4:G7").Value = "2016-08-18 08:32:44"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E4:E7").Value = "ht"
$wsZhCn.Range("H4:H7").Value = "2016-08-18 08:32:38"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E4:E7").Value = "ht"
$wsDeDe.Range("H4:H7").Value = "2016-08-18 08:32:44"
